# ------------------------------------------------------------------
# C1--C2-and-C3-PowerPoint.pptx commit replay
#   1) Re-style the 2-column data table on slide 16 with a different
#      (built-in) table style.
#   2) Re-colour the deck's theme from the "Integral" palette over to
#      the "Office" palette (the two a:clrScheme colour sets that ship
#      with this deck's theme parts).
# ------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
# Slide 16, shape 3 is the graphicFrame holding the "Total Outflow /
# Total Inflow" table. Table styles are not assignable through the
# .Style property directly - PowerPoint requires ApplyStyle(guid).
$tblShape = $p.Slides.Item(16).Shapes.Item(3)
$tblShape.Table.ApplyStyle("{B2C13479-82B7-4D7B-9DF2-336DD4B2FEAB}")

# --- 2. Theme colours -------------------------------------------------
# Swap the live theme's 12 scheme colours from the "Integral" values
# over to the "Office" values (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink - in that fixed order, matching ThemeColorScheme's indices
# 1-12). All slides/layouts share a single slide master, so editing
# through any one slide recolours the whole presentation.
function ConvertTo-OleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $themeColors.Item($i + 1).RGB = ConvertTo-OleColor $officeThemeColors[$i]
}
